$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '35.023.68'
$ws.Range("E2").Value = '  +0.45%  '
$ws.Range("D3").Value = '1.851.12'
$ws.Range("E3").Value = '  +2.54%  '
$ws.Range("E4").Value = '  +0.29%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '232.79'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.82%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.619'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.47%  '
$ws.Range("E7").Value = '  +0.31%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '40.75'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +4.56%  '
$ws.Range("E9").Value = '  +3.53%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0693'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.38%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0987'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.62%  '
$ws.Range("D12").Value = '2.120.27'
$ws.Range("E12").Value = '  +2.62%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '11.41'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +5.10%  '
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '1.855.39'
$ws.Range("E14").Value = '  +2.84%  '
$ws.Range("B15").Value = 'Polygon'
$ws.Range("C15").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.675'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.68%  '
$ws.Range("E16").Value = '  +2.92%  '
$ws.Range("D17").Value = '35.074.32'
$ws.Range("E17").Value = '  +0.74%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '70.13'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.67%  '
$ws.Range("E19").Value = '  +1.93%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '240.56'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.98%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.28'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +4.67%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.75'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.30%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.25'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.67%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '172.76'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.79%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.85'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.59%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.55'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.88%  '
$ws.Range("E28").Value = '  +4.09%  '
$ws.Range("E29").Value = '  +2.24%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0555'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.42%  '
$ws.Range("E32").Value = '  -0.22%  '
$ws.Range("E33").Value = '  +1.84%  '
$ws.Range("B34").Value = 'WEMIXToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.59'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +22.12%  '
$ws.Range("B35").Value = 'LidoDAOToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.96'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +12.05%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.757'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +11.03%  '
$ws.Range("E37").Value = '  +7.47%  '
$ws.Range("E38").Value = '  +13.15%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '90.59'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.39%  '
$ws.Range("D40").Value = '1.351.54'
$ws.Range("E40").Value = '  +3.75%  '
$ws.Range("E41").Value = '  +2.79%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '14.66'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.19%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.27'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.93%  '
$ws.Range("E44").Value = '  -1.67%  '
$ws.Range("E45").Value = '  +3.18%  '
$ws.Range("E46").Value = '  +4.51%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '6.33'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.45%  '
$ws.Range("D48").Value = '2.039.69'
$ws.Range("E48").Value = '  +2.73%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.42'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +19.61%  '
$ws.Range("E50").Value = '  +0.30%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0668'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.22%  '
